# Apply the "Results" sheet rework: abbreviate header labels, re-order the
# Syntax-Locked / Syntax-Free summary-stat columns (Mean/Median/Min/Max ahead
# of the Syntax-Free Count block) and resize the columns to fit the new
# shorter headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# --- Column widths -------------------------------------------------------
# COM ColumnWidth is in "characters"; the OOXML <col width> stored on save is
# ColumnWidth + 5/6 (the default 5px cell padding at Calibri 11's max digit
# width), so subtract 0.8333333333333334 to land exactly on the target width.
$ws.Columns.Item(2).ColumnWidth  = 19.166666666666668   # B -> 20
$ws.Columns.Item(3).ColumnWidth  = 19.166666666666668   # C -> 20
$ws.Columns.Item(6).ColumnWidth  = 8.166666666666666    # F -> 9
$ws.Columns.Item(7).ColumnWidth  = 19.166666666666668   # G -> 20
$ws.Columns.Item(8).ColumnWidth  = 9.166666666666666    # H -> 10
$ws.Columns.Item(9).ColumnWidth  = 6.166666666666667    # I -> 7
$ws.Columns.Item(10).ColumnWidth = 6.166666666666667    # J -> 7
$ws.Columns.Item(11).ColumnWidth = 8.166666666666666    # K -> 9
$ws.Columns.Item(12).ColumnWidth = 19.166666666666668   # L -> 20
$ws.Columns.Item(13).ColumnWidth = 9.166666666666666    # M -> 10
$ws.Columns.Item(14).ColumnWidth = 6.166666666666667    # N -> 7
$ws.Columns.Item(15).ColumnWidth = 6.166666666666667    # O -> 7

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("B1").Value = "SLStdDev"
$ws.Range("C1").Value = "SFStdDev"
$ws.Range("D1").Value = "PValue"
$ws.Range("E1").Value = "StdDevFactor"
$ws.Range("F1").Value = "SLCount"
$ws.Range("G1").Value = "SLMean"
$ws.Range("H1").Value = "SLMedian"
$ws.Range("I1").Value = "SLMin"
$ws.Range("J1").Value = "SLMax"
$ws.Range("K1").Value = "SFCount"
$ws.Range("L1").Value = "SFMean"
$ws.Range("M1").Value = "SFMedian"
$ws.Range("N1").Value = "SFMin"
$ws.Range("O1").Value = "SFMax"

# --- Data rows (2-5): re-map G:O into the new column order ---------------
# Old layout:  F=SLCount G=SFCount H=SLMean  I=SFMean  J=SLMedian
#              K=SFMedian L=SLMin  M=SFMin   N=SLMax   O=SFMax
# New layout:  F=SLCount G=SLMean  H=SLMedian I=SLMin  J=SLMax
#              K=SFCount L=SFMean  M=SFMedian N=SFMin  O=SFMax
for ($row = 2; $row -le 5; $row++) {
    $oldG = $ws.Cells.Item($row, 7).Value2   # SFCount
    $oldH = $ws.Cells.Item($row, 8).Value2   # SLMean
    $oldI = $ws.Cells.Item($row, 9).Value2   # SFMean
    $oldJ = $ws.Cells.Item($row, 10).Value2  # SLMedian
    $oldK = $ws.Cells.Item($row, 11).Value2  # SFMedian
    $oldL = $ws.Cells.Item($row, 12).Value2  # SLMin
    $oldM = $ws.Cells.Item($row, 13).Value2  # SFMin
    $oldN = $ws.Cells.Item($row, 14).Value2  # SLMax
    $oldO = $ws.Cells.Item($row, 15).Value2  # SFMax

    $ws.Cells.Item($row, 7).Value  = $oldH   # G: SLMean
    $ws.Cells.Item($row, 8).Value  = $oldJ   # H: SLMedian
    $ws.Cells.Item($row, 9).Value  = $oldL   # I: SLMin
    $ws.Cells.Item($row, 10).Value = $oldN   # J: SLMax
    $ws.Cells.Item($row, 11).Value = $oldG   # K: SFCount
    $ws.Cells.Item($row, 12).Value = $oldI   # L: SFMean
    $ws.Cells.Item($row, 13).Value = $oldK   # M: SFMedian
    $ws.Cells.Item($row, 14).Value = $oldM   # N: SFMin
    $ws.Cells.Item($row, 15).Value = $oldO   # O: SFMax
}

Write-Output "Results sheet header/column rework applied"
